$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[11, 35, 78, 123, 167, 246, 281, 366, 502, 609, 715, 755, 859]"
$ws.Range("D2").Value = "[17, 57, 136, 248, 356, 514, 766, 918, 1362, 1516, 2005, 2325, 2737]"
$ws.Range("E2").Value = "[2, 9, 21, 38, 54, 81, 125, 150, 228, 248, 318, 393, 455]"
$ws.Range("F2").Value = "[0.9333855344173657, 0.9546379362064165, 0.9035434077349098, 0.8053226360366034, 0.7249619150987372, 0.604319875930438, 0.5088548272825192, 0.3850235746730784, 0.24498917805420162, 0.21848481890473787, 0.11590220413675477, 0.08661139797050864, 0.05278898009786382]"

$ws.Range("C3").Value = "[12, 41, 74, 103, 144, 234, 271, 299, 367, 465, 519, 582, 737]"
$ws.Range("D3").Value = "[14, 56, 105, 158, 258, 407, 585, 704, 952, 1251, 1448, 1670, 2059]"
$ws.Range("E3").Value = "[2, 9, 18, 30, 46, 73, 103, 117, 160, 217, 249, 295, 366]"
$ws.Range("F3").Value = "[0.9918619131687658, 0.9569862094315925, 0.9092902045510506, 0.8470297967650031, 0.762929690000467, 0.6367971673481913, 0.5405187714385588, 0.49981810825897294, 0.397799615523347, 0.2642867237873717, 0.2130036236066909, 0.1456252319645504, 0.08849984908031076]"

$ws.Range("C4").Value = "[12, 41, 74, 109, 145, 213, 257, 337, 432, 537, 505, 718, 819]"
$ws.Range("D4").Value = "[14, 56, 105, 171, 256, 392, 538, 744, 931, 1242, 1374, 1873, 2217]"
$ws.Range("E4").Value = "[2, 9, 18, 30, 46, 70, 91, 132, 167, 228, 233, 339, 407]"
$ws.Range("F4").Value = "[0.9918619131687658, 0.9569862094315925, 0.9092902045510506, 0.8440360490253149, 0.7630999303838633, 0.6661534780293442, 0.5925685280846694, 0.44554571349759303, 0.3454116342015505, 0.2567005389107025, 0.23662510018397004, 0.12186712833352296, 0.06928845650696905]"

$ws.Range("C5").Value = "[12, 33, 61, 85, 119, 144, 193, 225, 255, 247, 320, 401, 371]"
$ws.Range("D5").Value = "[14, 47, 89, 175, 272, 392, 558, 727, 927, 1094, 1385, 1663, 1849]"
$ws.Range("E5").Value = "[2, 7, 14, 27, 43, 61, 90, 114, 144, 163, 211, 260, 284]"
$ws.Range("F5").Value = "[0.9341209303613215, 0.9595201151698435, 0.8370982389394316, 0.0, 0.7674899632460251, 0.6547107510348122, 0.39575427674543545, 0.4949941536587924, 0.3819051447804764, 0.2711603684251308, 0.23896966769282196, 0.0, 0.1451274985063912]"
